# update error compare UI
# The capture-image path column (J) used to be stored as a full relative
# path rooted at the use-case/site folder ("UC2.5/UC2.5_拠点A/..."); the UI
# now only needs the per-run path ("bdot.../....png"), so strip the
# redundant "UC2.5/UC2.5_拠点A/" prefix from every capimg cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "UC2.5/UC2.5_拠点A/"

$capimgCells = @("J5", "J7", "J8", "J9", "J10", "J11")
foreach ($addr in $capimgCells) {
    $cell = $ws.Range($addr)
    $text = $cell.Text
    if ($text -like "$prefix*") {
        $cell.Value = $text.Substring($prefix.Length)
    }
}

# Move the live selection to where the analyst is now working in the
# comparison sheet.
[void]$ws.Range("J15").Select()
